$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

# Insert a new row at position 65, shifting existing rows 65-162 down to 66-163
$ws.Rows.Item(65).Insert()

# Populate the newly inserted row 65 with the new record
$ws.Cells.Item(65, 1).Value = 5
$ws.Cells.Item(65, 2).Value = "Macroferia Regional de Talca"
$ws.Cells.Item(65, 3).Value = "Maule"
$ws.Cells.Item(65, 4).Value = 44571
$ws.Cells.Item(65, 5).Value = 7
$ws.Cells.Item(65, 6).Value = 100112021
$ws.Cells.Item(65, 7).Value = "Ají"
$ws.Cells.Item(65, 8).Value = "Americana (o)"
$ws.Cells.Item(65, 9).Value = "Primera"
$ws.Cells.Item(65, 10).Value = 150
$ws.Cells.Item(65, 11).Value = 18000
$ws.Cells.Item(65, 12).Value = 18000
$ws.Cells.Item(65, 13).Value = 18000
$ws.Cells.Item(65, 14).Value = "$/caja 14 kilos"
$ws.Cells.Item(65, 15).Value = "Región del Maule"
$ws.Cells.Item(65, 16).Value = 1286
$ws.Cells.Item(65, 17).Value = 14
$ws.Cells.Item(65, 18).Value = "Hortaliza"
